$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E14").Value = "menu background on main"
$ws.Range("A14").Value = "SOS main menu"
$ws.Range("B14").Value = "wav"
$ws.Range("C14").Value = "Jarryd"
$ws.Range("D14").Value = "n/a"

$ws.Range("E8").Select()
